$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column M mirrors column L's formatting (2020 -> 2021 data), so copy
# L3:L7 formats into M3:M7 first, then overwrite the values/content that
# differ from a straight copy.
$ws.Range("L3:L7").Copy()
$ws.Range("M3:M7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("M3").ClearContents()
$ws.Range("M4").Value = 2021
$ws.Range("M5").Value = 98
$ws.Range("M6").Value = 97
$ws.Range("M7").Value = 96

# The sheet was re-saved scrolled back to the top-left with the selection
# resting on A1 (instead of the stale topLeftCell="B1" / N13 selection it
# had before).
$ws.Range("A1").Select()
